# FutureShopTemplate-Valid.xlsx — add a missing "description" content-location
# entry to the Product_Tags block and duplicate the trailing "price" tag block
# (the sheet already repeats this html_tag / html_tag_attribute_name /
# html_tag_attribute_value pattern twice for the price; this adds a third
# occurrence right before the final content_location marker).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row above row 20 ("content_location") for a new
#        Item_Attribute_Name / description pair. Everything from the old
#        row 20 down shifts down by one row.
$ws.Rows("20:20").Insert()
$ws.Range("A20").Value = "Item_Attribute_Name"
$ws.Range("B20").Value = "description"

# --- 2. Insert a new row above what is now row 33 (the final
#        "content_location" row) for another Item_Attribute_Name / price
#        pair, duplicating the existing price-tag block a third time.
$ws.Rows("33:33").Insert()
$ws.Range("A33").Value = "Item_Attribute_Name"
$ws.Range("B33").Value = "price"

# --- 3. Update the view: scroll down a bit and move the selection.
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B23").Select()
